$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new data values (gang/spy strategy progress)
$ws.Range("D10").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("D22").Value = 1

# Force recalculation so the dependent formula cells (D14, E14, D25, E25, D27, E27) update
$excel.CalculateFullRebuild()

# Update the active selection to match the recorded cursor position
$ws.Range("M14").Select()
